# Update team-specific time-simulation matrix (Sheet1) with newly
# computed proportions for rows 2,3,6-13,15-19 per the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2105263157894737
$ws.Cells.Item(2, 3).Value = 0.5228070175438596
$ws.Cells.Item(2, 10).Value = 0.03859649122807018
$ws.Cells.Item(2, 16).Value = 0.1403508771929824
$ws.Cells.Item(2, 19).Value = 0.08771929824561403
$ws.Cells.Item(3, 2).Value = 0.02531645569620253
$ws.Cells.Item(3, 3).Value = 0.03164556962025317
$ws.Cells.Item(3, 10).Value = 0.0759493670886076
$ws.Cells.Item(3, 16).Value = 0.5949367088607594
$ws.Cells.Item(3, 19).Value = 0.2721518987341772
$ws.Cells.Item(6, 2).Value = 0.08994708994708994
$ws.Cells.Item(6, 4).Value = 0.01587301587301587
$ws.Cells.Item(6, 6).Value = 0.04232804232804233
$ws.Cells.Item(6, 10).Value = 0.1851851851851852
$ws.Cells.Item(6, 15).Value = 0.02116402116402116
$ws.Cells.Item(6, 17).Value = 0.1164021164021164
$ws.Cells.Item(6, 18).Value = 0.06349206349206349
$ws.Cells.Item(6, 19).Value = 0.4656084656084656
$ws.Cells.Item(7, 2).Value = 0.1465968586387434
$ws.Cells.Item(7, 4).Value = 0.02617801047120419
$ws.Cells.Item(7, 6).Value = 0.05759162303664921
$ws.Cells.Item(7, 10).Value = 0.1361256544502618
$ws.Cells.Item(7, 15).Value = 0.005235602094240838
$ws.Cells.Item(7, 17).Value = 0.1465968586387434
$ws.Cells.Item(7, 18).Value = 0.07329842931937172
$ws.Cells.Item(7, 19).Value = 0.4083769633507853
$ws.Cells.Item(8, 2).Value = 0.07142857142857142
$ws.Cells.Item(8, 4).Value = 0.00510204081632653
$ws.Cells.Item(8, 6).Value = 0.04336734693877551
$ws.Cells.Item(8, 10).Value = 0.09183673469387756
$ws.Cells.Item(8, 15).Value = 0.02040816326530612
$ws.Cells.Item(8, 17).Value = 0.2244897959183673
$ws.Cells.Item(8, 18).Value = 0.1147959183673469
$ws.Cells.Item(8, 19).Value = 0.4285714285714285
$ws.Cells.Item(9, 2).Value = 0.1304347826086956
$ws.Cells.Item(9, 4).Value = 0.005434782608695652
$ws.Cells.Item(9, 5).Value = 0.005434782608695652
$ws.Cells.Item(9, 6).Value = 0.07065217391304347
$ws.Cells.Item(9, 10).Value = 0.1413043478260869
$ws.Cells.Item(9, 15).Value = 0.02173913043478261
$ws.Cells.Item(9, 17).Value = 0.1684782608695652
$ws.Cells.Item(9, 18).Value = 0.09782608695652174
$ws.Cells.Item(9, 19).Value = 0.358695652173913
$ws.Cells.Item(10, 2).Value = 0.1019575856443719
$ws.Cells.Item(10, 4).Value = 0.02039151712887439
$ws.Cells.Item(10, 5).Value = 0.0008156606851549756
$ws.Cells.Item(10, 6).Value = 0.07340946166394779
$ws.Cells.Item(10, 10).Value = 0.1313213703099511
$ws.Cells.Item(10, 15).Value = 0.01957585644371941
$ws.Cells.Item(10, 17).Value = 0.2030995106035889
$ws.Cells.Item(10, 18).Value = 0.09053833605220228
$ws.Cells.Item(10, 19).Value = 0.3588907014681892
$ws.Cells.Item(11, 7).Value = 0.1066666666666667
$ws.Cells.Item(11, 10).Value = 0.1033333333333333
$ws.Cells.Item(11, 11).Value = 0.16
$ws.Cells.Item(11, 12).Value = 0.62
$ws.Cells.Item(11, 19).Value = 0.01
$ws.Cells.Item(12, 7).Value = 0.7424242424242424
$ws.Cells.Item(12, 10).Value = 0.1767676767676768
$ws.Cells.Item(12, 11).Value = 0.0101010101010101
$ws.Cells.Item(12, 12).Value = 0.02525252525252525
$ws.Cells.Item(12, 19).Value = 0.04545454545454546
$ws.Cells.Item(13, 7).Value = 0.6470588235294118
$ws.Cells.Item(13, 10).Value = 0.2941176470588235
$ws.Cells.Item(13, 19).Value = 0.05882352941176471
$ws.Cells.Item(15, 6).Value = 0.01984126984126984
$ws.Cells.Item(15, 8).Value = 0.123015873015873
$ws.Cells.Item(15, 9).Value = 0.09523809523809523
$ws.Cells.Item(15, 10).Value = 0.373015873015873
$ws.Cells.Item(15, 11).Value = 0.08333333333333333
$ws.Cells.Item(15, 13).Value = 0.007936507936507936
$ws.Cells.Item(15, 15).Value = 0.09126984126984126
$ws.Cells.Item(15, 19).Value = 0.2063492063492063
$ws.Cells.Item(16, 6).Value = 0.02027027027027027
$ws.Cells.Item(16, 8).Value = 0.1621621621621622
$ws.Cells.Item(16, 9).Value = 0.1081081081081081
$ws.Cells.Item(16, 10).Value = 0.4054054054054054
$ws.Cells.Item(16, 11).Value = 0.1216216216216216
$ws.Cells.Item(16, 13).Value = 0.02027027027027027
$ws.Cells.Item(16, 15).Value = 0.06081081081081081
$ws.Cells.Item(16, 19).Value = 0.1013513513513514
$ws.Cells.Item(17, 6).Value = 0.007211538461538462
$ws.Cells.Item(17, 8).Value = 0.1730769230769231
$ws.Cells.Item(17, 9).Value = 0.06971153846153846
$ws.Cells.Item(17, 10).Value = 0.4134615384615384
$ws.Cells.Item(17, 11).Value = 0.1129807692307692
$ws.Cells.Item(17, 13).Value = 0.02403846153846154
$ws.Cells.Item(17, 14).Value = 0.002403846153846154
$ws.Cells.Item(17, 15).Value = 0.08173076923076923
$ws.Cells.Item(17, 19).Value = 0.1153846153846154
$ws.Cells.Item(18, 6).Value = 0.01
$ws.Cells.Item(18, 8).Value = 0.19
$ws.Cells.Item(18, 9).Value = 0.07000000000000001
$ws.Cells.Item(18, 10).Value = 0.37
$ws.Cells.Item(18, 11).Value = 0.12
$ws.Cells.Item(18, 13).Value = 0.02
$ws.Cells.Item(18, 15).Value = 0.125
$ws.Cells.Item(18, 19).Value = 0.095
$ws.Cells.Item(19, 6).Value = 0.01357082273112807
$ws.Cells.Item(19, 8).Value = 0.1959287531806616
$ws.Cells.Item(19, 9).Value = 0.08566581849024597
$ws.Cells.Item(19, 10).Value = 0.3825275657336726
$ws.Cells.Item(19, 11).Value = 0.1153519932145886
$ws.Cells.Item(19, 13).Value = 0.01611535199321459
$ws.Cells.Item(19, 14).Value = 0.0008481764206955047
$ws.Cells.Item(19, 15).Value = 0.07718405428329092
$ws.Cells.Item(19, 19).Value = 0.1128074639525021
